$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Developer name
$ws.Range("C3").Value = "Jashanpreet Singh Sidhu"

# Preconditions column (E) for rows 7-11: set to "None"
$ws.Range("E7").Value = "None"
$ws.Range("E8").Value = "None"
$ws.Range("E9").Value = "None"
$ws.Range("E10").Value = "None"
$ws.Range("E11").Value = "None"

# Method Inputs column (F) for rows 12-16: set to "None"
$ws.Range("F12").Value = "None"
$ws.Range("F13").Value = "None"
$ws.Range("F14").Value = "None"
$ws.Range("F15").Value = "None"
$ws.Range("F16").Value = "None"

# Reflect final selection state seen in the saved file
$ws.Range("F17").Select()
